$wb = $excel.ActiveWorkbook

# Rename sheets (sheetId 1,3,4 change names; sheetId 2 unchanged)
$wb.Worksheets.Item(1).Name = "iCC389"
$wb.Worksheets.Item(3).Name = "iCC470"
$wb.Worksheets.Item(4).Name = "iCC651"

# Sheet 1 value updates
$ws = $wb.Worksheets.Item(1)
$ws.Range("B3").Value = 0.1334525431033358
$ws.Range("B6").Value = 0.1334525431033358
$ws.Range("B11").Value = 0.1332700813219138
$ws.Range("B12").Value = 0.1334525431033358
$ws.Range("B13").Value = 0.1329629536853232
$ws.Range("B17").Value = 0.1334525431033368
$ws.Range("B19").Value = 0.1334525431033355
$ws.Range("B20").Value = 0.1299261542360375
$ws.Range("B21").Value = 0.1315173027191182

# Sheet 2 value updates
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 0.9466355276173223
$ws.Range("B3").Value = 0.959935551660497
$ws.Range("B4").Value = 0.959935551660497
$ws.Range("B5").Value = 0.94871650757894
$ws.Range("B6").Value = 0.9526878130927015
$ws.Range("B7").Value = 0.9536983332068012
$ws.Range("B8").Value = 0.9585044486269945
$ws.Range("B10").Value = 0.9500172167264285
$ws.Range("B11").Value = 0.9494665968198088
$ws.Range("B12").Value = 0.9599355516604959
$ws.Range("B14").Value = 0.9301361339115282
$ws.Range("B15").Value = 0.9586036691791211
$ws.Range("B16").Value = 0.9444148413199831
$ws.Range("B18").Value = 0.959935551660502
$ws.Range("B19").Value = 0.9426647253939953
$ws.Range("B20").Value = 0.948421667482833
$ws.Range("B21").Value = 0.9468936649558083

# Sheet 3 value updates
$ws = $wb.Worksheets.Item(3)
$ws.Range("B13").Value = 0.7960639170575391
$ws.Range("B20").Value = 0.7677394430990008
$ws.Range("B21").Value = 0.7911654076940946

# Sheet 4 value updates
$ws = $wb.Worksheets.Item(4)
$ws.Range("B3").Value = 0.3152237567343069
$ws.Range("B6").Value = 0.3031150326334918
$ws.Range("B9").Value = 0.3166988434888303
$ws.Range("B11").Value = 0.3167956982730806
$ws.Range("B12").Value = 0.3095662209049149
$ws.Range("B13").Value = 0.3157073685763436
$ws.Range("B14").Value = 0.3127616688854793
$ws.Range("B15").Value = 0.3041745317278176
$ws.Range("B17").Value = 0.3130409577793544
$ws.Range("B18").Value = 0.3167956982730841
$ws.Range("B19").Value = 0.3023367435367467
$ws.Range("B20").Value = 0.3077682761953017
$ws.Range("B21").Value = 0.2951072661685817
